$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.471.87"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "2.106.22"
$ws.Range("E3").Value = "  +4.72%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'330.07"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").Value = "'0.5259"
$ws.Range("E7").Value = "  +2.45%  "

$ws.Range("D8").Value = "'0.4394"
$ws.Range("E8").Value = "  +3.22%  "

$ws.Range("D9").Value = "'0.08883"
$ws.Range("E9").Value = "  +1.60%  "

$ws.Range("D10").Value = "'48.36"
$ws.Range("E10").Value = "  +11.03%  "

$ws.Range("D11").Value = "'1.167"
$ws.Range("E11").Value = "  +2.95%  "

$ws.Range("D12").Value = "'24.88"
$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").Value = "2.111.92"
$ws.Range("E13").Value = "  +4.83%  "

$ws.Range("D14").Value = "'6.753"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").Value = "'7.773"
$ws.Range("E15").Value = "  +4.30%  "

$ws.Range("D16").Value = "'96.55"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").Value = "'0.00001130"
$ws.Range("E18").Value = "  +1.58%  "

$ws.Range("D19").Value = "'0.06639"
$ws.Range("E19").Value = "  +1.48%  "

$ws.Range("D20").Value = "'19.18"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'6.309"
$ws.Range("E22").Value = "  +1.69%  "

$ws.Range("D23").Value = "30.539.52"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "'12.27"
$ws.Range("E24").Value = "  +3.63%  "

$ws.Range("E25").Value = "  +4.22%  "

$ws.Range("D26").Value = "2.360.11"
$ws.Range("E26").Value = "  +4.80%  "

$ws.Range("D27").Value = "'22.48"
$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").Value = "'2.633"
$ws.Range("E28").Value = "  +8.05%  "

$ws.Range("D29").Value = "'161.89"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "'132.96"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").Value = "'1.232"
$ws.Range("E31").Value = "  +7.71%  "

$ws.Range("D32").Value = "'0.1074"
$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("D33").Value = "'1.677"
$ws.Range("E33").Value = "  +23.31%  "

$ws.Range("D34").Value = "'6.260"
$ws.Range("E34").Value = "  +2.80%  "

$ws.Range("D35").Value = "'3.898"
$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("D36").Value = "'10.20"
$ws.Range("E36").Value = "  +11.01%  "

$ws.Range("D37").Value = "'0.02594"
$ws.Range("E37").Value = "  +2.34%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06753"
$ws.Range("E38").Value = "  +1.21%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.515"
$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("D40").Value = "'12.74"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").Value = "'0.2285"
$ws.Range("E41").Value = "  +3.15%  "

$ws.Range("D42").Value = "'0.6923"
$ws.Range("E42").Value = "  +4.07%  "

$ws.Range("D43").Value = "'1.275"
$ws.Range("E43").Value = "  +2.82%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.13"
$ws.Range("E44").Value = "  +3.44%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6417"
$ws.Range("E46").Value = "  +3.85%  "

$ws.Range("D47").Value = "'2.225"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("D48").Value = "'3.632"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").Value = "'1.253"
$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").Value = "'1.216"
$ws.Range("E50").Value = "  +9.99%  "

$ws.Range("D51").Value = "'82.74"
